$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: clear K6
$ws.Range("K6").Value = ""

# Row 12: update D12, clear K12
$ws.Range("D12").Value = "Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = ""

# Row 13: update E13, set F13, clear K13
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 11.380 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 795 TL | Azami 4.005 TL"
$ws.Range("K13").Value = ""

# Row 14: set F14, clear K14
$ws.Range("F14").Value = "2.785,72 TL - 12.380,95 TL"
$ws.Range("K14").Value = ""
